$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 40000
$ws.Range("J3").Value = 40000
$ws.Range("L3").Value = 40000
$ws.Range("N3").Value = -40228
$ws.Range("H6").Value = 1002
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").Value = $null
$ws.Range("H17").Value = 1537.4166
$ws.Range("J17").Value = 1651.1
$ws.Range("L17").Value = 4953.299999999999
$ws.Range("N17").Value = -5289.299999999999
$ws.Range("H19").Value = 1975
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1975
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = $null
$ws.Range("M19").Value = 1975
$ws.Range("N19").Value = -2325
$ws.Range("H33").Value = 143.8421
$ws.Range("I33").Value = 142.6875
$ws.Range("K33").Value = 142.6875
$ws.Range("M33").Value = 86.3125
$ws.Range("H53").Value = 226.05556
$ws.Range("J53").Value = 206.22223
$ws.Range("L53").Value = 206.22223
$ws.Range("N53").Value = -1480.22223
$ws.Range("H92").Value = 675.4
$ws.Range("I92").Value = 706.6842
$ws.Range("J92").Value = 576.3333
$ws.Range("K92").Value = 706.6842
$ws.Range("L92").Value = 576.3333
$ws.Range("M92").Value = 541.3158
$ws.Range("N92").Value = -3072.3333
$ws.Range("H102").Value = 40000
$ws.Range("J102").Value = 40000
$ws.Range("L102").Value = 40000
$ws.Range("N102").Value = -46490
$ws.Range("H116").Value = 20523.95
$ws.Range("I116").Value = 17099.834
$ws.Range("K116").Value = 17099.834
$ws.Range("M116").Value = -13657.834
$ws.Range("H129").Value = 2641.2
$ws.Range("I129").Value = 428.8
$ws.Range("K129").Value = 1286.4
$ws.Range("M129").Value = 3713.6
$ws.Range("H132").Value = 8116.6665
$ws.Range("I132").Value = 8116.6665
$ws.Range("K132").Value = 24349.9995
$ws.Range("M132").Value = -21819.9995
$ws.Range("H137").Value = 6151.364
$ws.Range("I137").Value = 1917.3846
$ws.Range("J137").Value = 12267.111
$ws.Range("K137").Value = 5752.1538
$ws.Range("L137").Value = 36801.333
$ws.Range("M137").Value = -3202.1538
$ws.Range("N137").Value = -41901.333
$ws.Range("H138").Value = 4134.4287
$ws.Range("I138").Value = 5212.2
$ws.Range("K138").Value = 15636.6
$ws.Range("M138").Value = -10496.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 210427.61
$ws.Range("J32").Value = 39916.332
$ws.Range("L32").Value = 39916.332
$ws.Range("N32").Value = -40490.332
$ws.Range("H45").Value = 5000
$ws.Range("I45").Value = 2000
$ws.Range("K45").Value = 2000
$ws.Range("M45").Value = -1623
$ws.Range("H61").Value = 8999.5
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 8999.5
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = $null
$ws.Range("M61").Value = 8999.5
$ws.Range("N61").Value = -9423.5
$ws.Range("H74").Value = 12945.952
$ws.Range("I74").Value = 7473.75
$ws.Range("K74").Value = 7473.75
$ws.Range("M74").Value = -6599.75
$ws.Range("H77").Value = 12945.952
$ws.Range("I77").Value = 7473.75
$ws.Range("K77").Value = 37368.75
$ws.Range("M77").Value = -33000.75
$ws.Range("H107").Value = 40228
$ws.Range("J107").Value = 40228
$ws.Range("L107").Value = 40228
$ws.Range("N107").Value = -47908
$ws.Range("H122").Value = 2334
$ws.Range("I122").Value = 2103.3845
$ws.Range("K122").Value = 6310.1535
$ws.Range("M122").Value = -3860.1535
$ws.Range("H128").Value = 30000
$ws.Range("J128").Value = 30000
$ws.Range("L128").Value = 30000
$ws.Range("N128").Value = -39960
$ws.Range("H132").Value = 6876.9165
$ws.Range("I132").Value = 4453.8335
$ws.Range("J132").Value = 9300
$ws.Range("K132").Value = 13361.5005
$ws.Range("L132").Value = 27900
$ws.Range("M132").Value = -10831.5005
$ws.Range("N132").Value = -32960
$ws.Range("H136").Value = 8999.5
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 8999.5
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = $null
$ws.Range("M136").Value = 26998.5
$ws.Range("N136").Value = -32098.5
$ws.Range("H138").Value = 80194.5
$ws.Range("I138").Value = 80390
$ws.Range("J138").Value = 79999
$ws.Range("K138").Value = 80390
$ws.Range("L138").Value = 79999
$ws.Range("M138").Value = -75250
$ws.Range("N138").Value = -90279

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 95999.5
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").Value = $null
$ws.Range("H58").Value = 57042.25
$ws.Range("J58").Value = 64153.332
$ws.Range("L58").Value = 64153.332
$ws.Range("N58").Value = -64741.332
$ws.Range("H80").Value = 899.63635
$ws.Range("I80").Value = 806.9286
$ws.Range("J80").Value = 1061.875
$ws.Range("K80").Value = 806.9286
$ws.Range("L80").Value = 1061.875
$ws.Range("M80").Value = 191.0714
$ws.Range("N80").Value = -3057.875
$ws.Range("H83").Value = 899.63635
$ws.Range("I83").Value = 806.9286
$ws.Range("J83").Value = 1061.875
$ws.Range("K83").Value = 4034.643
$ws.Range("L83").Value = 5309.375
$ws.Range("M83").Value = 957.357
$ws.Range("N83").Value = -15293.375
$ws.Range("H94").Value = 2408.8696
$ws.Range("I94").Value = 1694.7
$ws.Range("K94").Value = 1694.7
$ws.Range("M94").Value = -1243.7
$ws.Range("H99").Value = 9238.857
$ws.Range("I99").Value = 12844.333
$ws.Range("K99").Value = 12844.333
$ws.Range("M99").Value = -11346.333
$ws.Range("H105").Value = 2127.1428
$ws.Range("I105").Value = 1864.6538
$ws.Range("J105").Value = 2885.4443
$ws.Range("K105").Value = 1864.6538
$ws.Range("L105").Value = 2885.4443
$ws.Range("M105").Value = -117.6538
$ws.Range("N105").Value = -6379.4443
$ws.Range("H107").Value = 1137.45
$ws.Range("I107").Value = 1132.2941
$ws.Range("J107").Value = 1166.6666
$ws.Range("K107").Value = 1132.2941
$ws.Range("L107").Value = 1166.6666
$ws.Range("M107").Value = 787.7058999999999
$ws.Range("N107").Value = -5006.6666
$ws.Range("H134").Value = 6791.364
$ws.Range("I134").Value = 3676.64
$ws.Range("J134").Value = 16524.875
$ws.Range("K134").Value = 11029.92
$ws.Range("L134").Value = 49574.625
$ws.Range("M134").Value = -8494.92
$ws.Range("N134").Value = -54644.625
$ws.Range("H136").Value = 95999.5
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = $null
$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 79000.39999999999
$ws.Range("I16").Value = 18402.375
$ws.Range("J16").Value = 148255.28
$ws.Range("K16").Value = 18402.375
$ws.Range("L16").Value = 148255.28
$ws.Range("M16").Value = -18115.375
$ws.Range("N16").Value = -148829.28
$ws.Range("H19").Value = 771.5714
$ws.Range("I19").Value = 741.0909
$ws.Range("K19").Value = 741.0909
$ws.Range("M19").Value = -571.0909
$ws.Range("H22").Value = 2255.625
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 2449.2856
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 2449.2856
$ws.Range("M22").Value = -550
$ws.Range("N22").Value = -3149.2856
$ws.Range("H24").Value = 771.5714
$ws.Range("I24").Value = 741.0909
$ws.Range("K24").Value = 741.0909
$ws.Range("M24").Value = -571.0909
$ws.Range("H31").Value = 2772.889
$ws.Range("I31").Value = 4660
$ws.Range("J31").Value = 1829.3334
$ws.Range("K31").Value = 4660
$ws.Range("L31").Value = 1829.3334
$ws.Range("M31").Value = -4365
$ws.Range("N31").Value = -2419.3334
$ws.Range("H34").Value = 2772.889
$ws.Range("I34").Value = 4660
$ws.Range("J34").Value = 1829.3334
$ws.Range("K34").Value = 4660
$ws.Range("L34").Value = 1829.3334
$ws.Range("M34").Value = -4458
$ws.Range("N34").Value = -2233.3334
$ws.Range("H58").Value = 4724.1665
$ws.Range("I58").Value = 3344.72
$ws.Range("K58").Value = 3344.72
$ws.Range("M58").Value = -3141.72
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = $null
$ws.Range("N92").Value = 0
$ws.Range("H99").Value = 12174.318
$ws.Range("I99").Value = 18407
$ws.Range("K99").Value = 18407
$ws.Range("M99").Value = -16909
$ws.Range("H105").Value = 11682.556
$ws.Range("I105").Value = 13042.875
$ws.Range("K105").Value = 13042.875
$ws.Range("M105").Value = -11295.875
$ws.Range("H107").Value = 1068.3846
$ws.Range("I107").Value = 723.875
$ws.Range("K107").Value = 723.875
$ws.Range("M107").Value = 1196.125
$ws.Range("H113").Value = 79000.39999999999
$ws.Range("I113").Value = 18402.375
$ws.Range("J113").Value = 148255.28
$ws.Range("K113").Value = 18402.375
$ws.Range("L113").Value = 148255.28
$ws.Range("M113").Value = -16232.375
$ws.Range("N113").Value = -152595.28
$ws.Range("H122").Value = 7592.9556
$ws.Range("I122").Value = 2075.6487
$ws.Range("J122").Value = 33110.5
$ws.Range("K122").Value = 6226.946100000001
$ws.Range("L122").Value = 99331.5
$ws.Range("M122").Value = -3776.946100000001
$ws.Range("N122").Value = -104231.5
$ws.Range("H126").Value = 12174.318
$ws.Range("I126").Value = 18407
$ws.Range("K126").Value = 55221
$ws.Range("M126").Value = -52751
$ws.Range("H132").Value = 3615.3225
$ws.Range("I132").Value = 3491.9644
$ws.Range("J132").Value = 4766.6665
$ws.Range("K132").Value = 10475.8932
$ws.Range("L132").Value = 14299.9995
$ws.Range("M132").Value = -7945.893199999999
$ws.Range("N132").Value = -19359.9995
$ws.Range("H134").Value = 5073.8
$ws.Range("I134").Value = 5642.25
$ws.Range("J134").Value = 2800
$ws.Range("K134").Value = 16926.75
$ws.Range("L134").Value = 8400
$ws.Range("M134").Value = -14391.75
$ws.Range("N134").Value = -13470
$ws.Range("H136").Value = 4724.1665
$ws.Range("I136").Value = 3344.72
$ws.Range("K136").Value = 10034.16
$ws.Range("M136").Value = -7484.16

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3000136
$ws.Range("I4").Value = 3230877.2
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 9692631.600000001
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -9692519.600000001
$ws.Range("N4").Value = -1724
$ws.Range("H5").Value = 1031.2963
$ws.Range("I5").Value = 849
$ws.Range("J5").Value = 1341.2
$ws.Range("K5").Value = 2547
$ws.Range("L5").Value = 4023.6
$ws.Range("M5").Value = -2435
$ws.Range("N5").Value = -4247.6
$ws.Range("H7").Value = 30
$ws.Range("I7").Value = 12.5
$ws.Range("K7").Value = 37.5
$ws.Range("M7").Value = 74.5
$ws.Range("H10").Value = 683.44446
$ws.Range("I10").Value = 164.14285
$ws.Range("J10").Value = 2501
$ws.Range("K10").Value = 492.42855
$ws.Range("L10").Value = 7503
$ws.Range("M10").Value = -353.42855
$ws.Range("N10").Value = -7781
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = $null
$ws.Range("H70").Value = 676.8333
$ws.Range("I70").Value = 676.8333
$ws.Range("K70").Value = 2030.4999
$ws.Range("M70").Value = -1715.4999
$ws.Range("H73").Value = 676.8333
$ws.Range("I73").Value = 676.8333
$ws.Range("K73").Value = 2030.4999
$ws.Range("M73").Value = -938.4999
$ws.Range("H92").Value = 784.78125
$ws.Range("I92").Value = 485.61905
$ws.Range("J92").Value = 1355.909
$ws.Range("K92").Value = 1456.85715
$ws.Range("L92").Value = 4067.727
$ws.Range("M92").Value = -208.85715
$ws.Range("N92").Value = -6563.727000000001
$ws.Range("H114").Value = 1047
$ws.Range("J114").Value = 2313
$ws.Range("L114").Value = 6939
$ws.Range("N114").Value = -13447
$ws.Range("H122").Value = 1242297.5
$ws.Range("I122").Value = 5376507
$ws.Range("J122").Value = 2034.7
$ws.Range("K122").Value = 48388563
$ws.Range("L122").Value = 18312.3
$ws.Range("M122").Value = -48386113
$ws.Range("N122").Value = -23212.3
$ws.Range("H131").Value = 1489.5128
$ws.Range("J131").Value = 1622.742
$ws.Range("L131").Value = 4868.226
$ws.Range("N131").Value = -14948.226
$ws.Range("H135").Value = 1031.2963
$ws.Range("I135").Value = 849
$ws.Range("J135").Value = 1341.2
$ws.Range("K135").Value = 7641
$ws.Range("L135").Value = 12070.8
$ws.Range("M135").Value = -5106
$ws.Range("N135").Value = -17140.8
$ws.Range("H136").Value = 6727.3125
$ws.Range("I136").Value = 7013.9
$ws.Range("J136").Value = 6249.6665
$ws.Range("K136").Value = 21041.7
$ws.Range("L136").Value = 18748.9995
$ws.Range("M136").Value = -15941.7
$ws.Range("N136").Value = -28948.9995
$ws.Range("H137").Value = 3276.4167
$ws.Range("I137").Value = 1911.75
$ws.Range("J137").Value = 6005.75
$ws.Range("K137").Value = 5735.25
$ws.Range("L137").Value = 18017.25
$ws.Range("M137").Value = -635.25
$ws.Range("N137").Value = -28217.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = $null
$ws.Range("N53").Value = 0
$ws.Range("H70").Value = 40227.668
$ws.Range("I70").Value = 44466.625
$ws.Range("J70").Value = 31749.75
$ws.Range("K70").Value = 44466.625
$ws.Range("L70").Value = 31749.75
$ws.Range("M70").Value = -44196.625
$ws.Range("N70").Value = -32289.75
$ws.Range("H73").Value = 40227.668
$ws.Range("I73").Value = 44466.625
$ws.Range("J73").Value = 31749.75
$ws.Range("K73").Value = 44466.625
$ws.Range("L73").Value = 31749.75
$ws.Range("M73").Value = -43530.625
$ws.Range("N73").Value = -33621.75
$ws.Range("H80").Value = 3044.3572
$ws.Range("I80").Value = 2820.4285
$ws.Range("J80").Value = 3268.2856
$ws.Range("K80").Value = 2820.4285
$ws.Range("L80").Value = 3268.2856
$ws.Range("M80").Value = -1822.4285
$ws.Range("N80").Value = -5264.2856
$ws.Range("H83").Value = 3044.3572
$ws.Range("I83").Value = 2820.4285
$ws.Range("J83").Value = 3268.2856
$ws.Range("K83").Value = 14102.1425
$ws.Range("L83").Value = 16341.428
$ws.Range("M83").Value = -9110.1425
$ws.Range("N83").Value = -26325.428
$ws.Range("H102").Value = 2551.3333
$ws.Range("I102").Value = 2756.6667
$ws.Range("J102").Value = 1730
$ws.Range("K102").Value = 2756.6667
$ws.Range("L102").Value = 1730
$ws.Range("M102").Value = -1134.6667
$ws.Range("N102").Value = -4974
$ws.Range("H122").Value = 7286.5713
$ws.Range("I122").Value = 11002
$ws.Range("K122").Value = 33006
$ws.Range("M122").Value = -30556
$ws.Range("H132").Value = 54506
$ws.Range("I132").Value = 54506
$ws.Range("K132").Value = 163518
$ws.Range("M132").Value = -160988
$ws.Range("H136").Value = 63335.637
$ws.Range("J136").Value = 63335.637
$ws.Range("L136").Value = 190006.911
$ws.Range("N136").Value = -195106.911
$ws.Range("H138").Value = 69666.336
$ws.Range("I138").Value = 10000
$ws.Range("K138").Value = 10000
$ws.Range("M138").Value = -4860

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 635.4286
$ws.Range("I16").Value = 592
$ws.Range("K16").Value = 592
$ws.Range("M16").Value = -422
$ws.Range("H40").Value = 2142.8667
$ws.Range("I40").Value = 2165.077
$ws.Range("J40").Value = 1998.5
$ws.Range("K40").Value = 2165.077
$ws.Range("L40").Value = 1998.5
$ws.Range("M40").Value = -2029.077
$ws.Range("N40").Value = -2270.5
$ws.Range("H43").Value = 34000
$ws.Range("J43").Value = 34000
$ws.Range("L43").Value = 34000
$ws.Range("N43").Value = -34386
$ws.Range("H93").Value = 4656.3335
$ws.Range("I93").Value = 2000
$ws.Range("J93").Value = 5984.5
$ws.Range("K93").Value = 2000
$ws.Range("L93").Value = 5984.5
$ws.Range("M93").Value = -752
$ws.Range("N93").Value = -8480.5
$ws.Range("H114").Value = 55000
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").Value = $null
$ws.Range("H125").Value = 71220.664
$ws.Range("J125").Value = 71220.664
$ws.Range("L125").Value = 71220.664
$ws.Range("N125").Value = -81060.664
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = $null
$ws.Range("N128").Value = 0

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 10000
$ws.Range("K26").Value = 10000
$ws.Range("M26").Value = -9707
$ws.Range("H81").Value = 71503480
$ws.Range("J81").Value = 166837620
$ws.Range("L81").Value = 333675240
$ws.Range("N81").Value = -333677362
$ws.Range("H84").Value = 71503480
$ws.Range("J84").Value = 166837620
$ws.Range("L84").Value = 1668376200
$ws.Range("N84").Value = -1668386808
$ws.Range("H100").Value = 13889829
$ws.Range("I100").Value = 1259.8
$ws.Range("K100").Value = 2519.6
$ws.Range("M100").Value = -1978.6
$ws.Range("H107").Value = 26317916
$ws.Range("I107").Value = 37037932
$ws.Range("J107").Value = 5152.727
$ws.Range("K107").Value = 111113796
$ws.Range("L107").Value = 15458.181
$ws.Range("M107").Value = -111111876
$ws.Range("N107").Value = -19298.181
$ws.Range("H109").Value = 35000
$ws.Range("J109").Value = 35000
$ws.Range("L109").Value = 35000
$ws.Range("N109").Value = -37774
$ws.Range("H118").Value = 50000
$ws.Range("J118").Value = 50000
$ws.Range("L118").Value = 50000
$ws.Range("N118").Value = -53314
$ws.Range("H122").Value = 28850.559
$ws.Range("I122").Value = 2445.6765
$ws.Range("K122").Value = 7337.029500000001
$ws.Range("M122").Value = -4887.029500000001
$ws.Range("H132").Value = 2060.1072
$ws.Range("I132").Value = 1783.409
$ws.Range("J132").Value = 3074.6667
$ws.Range("K132").Value = 5350.227000000001
$ws.Range("L132").Value = 9224.000100000001
$ws.Range("M132").Value = -2820.227000000001
$ws.Range("N132").Value = -14284.0001
$ws.Range("H136").Value = 2747.6667
$ws.Range("I136").Value = 2169.375
$ws.Range("J136").Value = 4598.2
$ws.Range("K136").Value = 6508.125
$ws.Range("L136").Value = 13794.6
$ws.Range("M136").Value = -3958.125
$ws.Range("N136").Value = -18894.6
